$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_vals data (filtered save games) for rows 2-25, columns B-G
$ws.Cells.Item(2, 2).Value = 3.230985683306322
$ws.Cells.Item(2, 3).Value = 1.667794583268128
$ws.Cells.Item(2, 4).Value = 0.8054896365839992
$ws.Cells.Item(2, 5).Value = 8.660232485948974
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 14.36450238910742
$ws.Cells.Item(3, 2).Value = 3.230985683306322
$ws.Cells.Item(3, 3).Value = 1.667794583268128
$ws.Cells.Item(3, 4).Value = 0.8054896365839992
$ws.Cells.Item(3, 5).Value = 0.496779210170732
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 6.201049113329182
$ws.Cells.Item(4, 2).Value = 0.3048080303191223
$ws.Cells.Item(4, 3).Value = 0.04240448674262143
$ws.Cells.Item(4, 4).Value = 0.1575252929769615
$ws.Cells.Item(4, 5).Value = 0.496779210170732
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.001517020209437
$ws.Cells.Item(5, 2).Value = 3.230985683306322
$ws.Cells.Item(5, 3).Value = 1.667794583268128
$ws.Cells.Item(5, 4).Value = 0.1575252929769615
$ws.Cells.Item(5, 5).Value = 0.496779210170732
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 5.553084769722144
$ws.Cells.Item(6, 2).Value = 1.459612070389937
$ws.Cells.Item(6, 3).Value = 0.3127903958511391
$ws.Cells.Item(6, 4).Value = 0.1575252929769615
$ws.Cells.Item(6, 5).Value = 0.496779210170732
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 2.42670696938877
$ws.Cells.Item(7, 2).Value = 0.01514828764759746
$ws.Cells.Item(7, 3).Value = 0.04240448674262143
$ws.Cells.Item(7, 4).Value = 0.8054896365839992
$ws.Cells.Item(7, 5).Value = 0.496779210170732
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.35982162114495
$ws.Cells.Item(8, 2).Value = 1.459612070389937
$ws.Cells.Item(8, 3).Value = 1.667794583268128
$ws.Cells.Item(8, 4).Value = 3.900430680208489
$ws.Cells.Item(8, 5).Value = 0.496779210170732
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 7.524616544037286
$ws.Cells.Item(9, 2).Value = 1.459612070389937
$ws.Cells.Item(9, 3).Value = 1.667794583268128
$ws.Cells.Item(9, 4).Value = 3.900430680208489
$ws.Cells.Item(9, 5).Value = 0.496779210170732
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 7.524616544037286
$ws.Cells.Item(10, 2).Value = 1.459612070389937
$ws.Cells.Item(10, 3).Value = 1.667794583268128
$ws.Cells.Item(10, 4).Value = 3.900430680208489
$ws.Cells.Item(10, 5).Value = 8.660232485948974
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 15.68806981981553
$ws.Cells.Item(11, 2).Value = 0.04763786555579896
$ws.Cells.Item(11, 3).Value = 0.3127903958511391
$ws.Cells.Item(11, 4).Value = 3.900430680208489
$ws.Cells.Item(11, 5).Value = 0.496779210170732
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4.757638151786159
$ws.Cells.Item(12, 2).Value = 1.459612070389937
$ws.Cells.Item(12, 3).Value = 1.667794583268128
$ws.Cells.Item(12, 4).Value = 0.8054896365839992
$ws.Cells.Item(12, 5).Value = 0.496779210170732
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 4.429675500412797
$ws.Cells.Item(13, 2).Value = 3.230985683306322
$ws.Cells.Item(13, 3).Value = 1.667794583268128
$ws.Cells.Item(13, 4).Value = 0.8054896365839992
$ws.Cells.Item(13, 5).Value = 0.496779210170732
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 6.201049113329182
$ws.Cells.Item(14, 2).Value = 3.230985683306322
$ws.Cells.Item(14, 3).Value = 1.667794583268128
$ws.Cells.Item(14, 4).Value = 0.1575252929769615
$ws.Cells.Item(14, 5).Value = 0.496779210170732
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 5.553084769722144
$ws.Cells.Item(15, 2).Value = 0.127881588408715
$ws.Cells.Item(15, 3).Value = 1.667794583268128
$ws.Cells.Item(15, 4).Value = 0.8054896365839992
$ws.Cells.Item(15, 5).Value = 8.660232485948974
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 11.26139829420982
$ws.Cells.Item(16, 2).Value = 0.127881588408715
$ws.Cells.Item(16, 3).Value = 1.667794583268128
$ws.Cells.Item(16, 4).Value = 0.8054896365839992
$ws.Cells.Item(16, 5).Value = 0.496779210170732
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 3.097945018431574
$ws.Cells.Item(17, 2).Value = 0.127881588408715
$ws.Cells.Item(17, 3).Value = 0.3127903958511391
$ws.Cells.Item(17, 4).Value = 26.21740644021617
$ws.Cells.Item(17, 5).Value = 8.660232485948974
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 35.318310910425
$ws.Cells.Item(18, 2).Value = 3.230985683306322
$ws.Cells.Item(18, 3).Value = 1.667794583268128
$ws.Cells.Item(18, 4).Value = 0.8054896365839992
$ws.Cells.Item(18, 5).Value = 0.496779210170732
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 6.201049113329182
$ws.Cells.Item(19, 2).Value = 3.230985683306322
$ws.Cells.Item(19, 3).Value = 1.667794583268128
$ws.Cells.Item(19, 4).Value = 0.8054896365839992
$ws.Cells.Item(19, 5).Value = 0.496779210170732
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 6.201049113329182
$ws.Cells.Item(20, 2).Value = 0.6753301551942219
$ws.Cells.Item(20, 3).Value = 1.667794583268128
$ws.Cells.Item(20, 4).Value = 0.8054896365839992
$ws.Cells.Item(20, 5).Value = 0.496779210170732
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 3.645393585217082
$ws.Cells.Item(21, 2).Value = 3.230985683306322
$ws.Cells.Item(21, 3).Value = 1.667794583268128
$ws.Cells.Item(21, 4).Value = 0.8054896365839992
$ws.Cells.Item(21, 5).Value = 0.496779210170732
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 6.201049113329182
$ws.Cells.Item(22, 2).Value = 3.230985683306322
$ws.Cells.Item(22, 3).Value = 1.667794583268128
$ws.Cells.Item(22, 4).Value = 3.900430680208489
$ws.Cells.Item(22, 5).Value = 0.496779210170732
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = 9.295990156953671
$ws.Cells.Item(23, 2).Value = 3.230985683306322
$ws.Cells.Item(23, 3).Value = 1.667794583268128
$ws.Cells.Item(23, 4).Value = 0.1575252929769615
$ws.Cells.Item(23, 5).Value = 0.496779210170732
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 5.553084769722144
$ws.Cells.Item(24, 2).Value = 1.459612070389937
$ws.Cells.Item(24, 3).Value = 1.667794583268128
$ws.Cells.Item(24, 4).Value = 0.8054896365839992
$ws.Cells.Item(24, 5).Value = 0.496779210170732
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = 4.429675500412797
$ws.Cells.Item(25, 2).Value = 3.230985683306322
$ws.Cells.Item(25, 3).Value = 1.667794583268128
$ws.Cells.Item(25, 4).Value = 0.8054896365839992
$ws.Cells.Item(25, 5).Value = 0.496779210170732
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(25, 7).Value = 6.201049113329182

Write-Output "Updated B2:G25 with regenerated s_vals data"
